$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores values as text (e.g. "51.367.98", "1.00"),
# not numbers. Force text format on the D cells we are about to rewrite so
# Excel does not reinterpret values like "1.00" or "384.61" as numbers.
$priceRows = @(2,3,4,5,6,7,10,13,14,15,16,17,18,19,21,23,24,26,27,28,29,33,35,36,39,44,46,47,49,50,51)
foreach ($r in $priceRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

$ws.Range("D2").Value = "51.385.41"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "3.043.33"
$ws.Range("E3").Value = "  +2.54%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "384.61"
$ws.Range("E5").Value = "  +1.12%  "
$ws.Range("D6").Value = "102.77"
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("D7").Value = "0.543"
$ws.Range("E7").Value = "  -0.50%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -1.39%  "
$ws.Range("D10").Value = "36.70"
$ws.Range("E10").Value = "  +0.30%  "
$ws.Range("E11").Value = "  +0.06%  "
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("D13").Value = "3.526.90"
$ws.Range("E13").Value = "  +2.85%  "
$ws.Range("D14").Value = "18.56"
$ws.Range("E14").Value = "  +1.17%  "
$ws.Range("D15").Value = "7.74"
$ws.Range("E15").Value = "  -0.71%  "
$ws.Range("D16").Value = "3.043.42"
$ws.Range("E16").Value = "  +3.09%  "
$ws.Range("D17").Value = "0.969"
$ws.Range("E17").Value = "  -2.51%  "
$ws.Range("D18").Value = "10.54"
$ws.Range("E18").Value = "  -5.07%  "
$ws.Range("D19").Value = "51.465.50"
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("D21").Value = "12.38"
$ws.Range("E21").Value = "  -1.52%  "
$ws.Range("E22").Value = "  +0.28%  "
$ws.Range("D23").Value = "70.09"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").Value = "267.67"
$ws.Range("E24").Value = "  +0.41%  "
$ws.Range("E25").Value = "  -1.56%  "
$ws.Range("D26").Value = "8.18"
$ws.Range("E26").Value = "  +4.08%  "
$ws.Range("D27").Value = "26.91"
$ws.Range("E27").Value = "  +3.49%  "
$ws.Range("D28").Value = "0.170"
$ws.Range("D29").Value = "7.25"
$ws.Range("E29").Value = "  -3.45%  "
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("E31").Value = "  -1.98%  "
$ws.Range("E32").Value = "  -0.63%  "
$ws.Range("D33").Value = "34.78"
$ws.Range("E33").Value = "  +0.22%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").Value = "50.33"
$ws.Range("E35").Value = "  -2.13%  "
$ws.Range("D36").Value = "0.0446"
$ws.Range("E36").Value = "  +2.12%  "
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("E38").Value = "  +2.74%  "
$ws.Range("D39").Value = "0.290"
$ws.Range("E39").Value = "  +7.64%  "
$ws.Range("E40").Value = "  +1.91%  "
$ws.Range("E41").Value = "  +1.12%  "
$ws.Range("E42").Value = "  -0.86%  "
$ws.Range("E43").Value = "  -0.22%  "
$ws.Range("D44").Value = "124.78"
$ws.Range("E44").Value = "  +0.41%  "
$ws.Range("E45").Value = "  +3.43%  "
$ws.Range("D46").Value = "21.74"
$ws.Range("E46").Value = "  +0.80%  "
$ws.Range("D47").Value = "2.09"
$ws.Range("E47").Value = "  +3.01%  "
$ws.Range("E48").Value = "  +1.97%  "
$ws.Range("D49").Value = "2.026.02"
$ws.Range("E49").Value = "  -0.32%  "
$ws.Range("D50").Value = "3.343.05"
$ws.Range("E50").Value = "  +2.67%  "
$ws.Range("D51").Value = "0.0318"
$ws.Range("E51").Value = "  -3.36%  "
